$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25/26 swap: Monero <-> LidoDAOToken (B,C,D,E columns; A stays as rank index)
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.216"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'142.26"
$ws.Range("E26").Value = "  -2.28%  "

$ws.Range("D2").Value = "26.613.01"
$ws.Range("D3").Value = "1.842.95"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D5").Value = "'259.25"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5274"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").Value = "'0.3150"
$ws.Range("E8").Value = "  -3.79%  "
$ws.Range("D9").Value = "'0.06804"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'18.76"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").Value = "'0.7810"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "'0.07754"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "1.850.10"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "'88.02"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "'5.012"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'13.86"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'0.000007934"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "26.640.94"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "2.072.88"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "'4.613"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'5.978"
$ws.Range("D24").Value = "'9.316"
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").Value = "'17.00"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "'110.87"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "'4.189"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'0.08734"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'4.080"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'0.04877"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "'0.7284"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("D36").Value = "'2.857"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'3.091"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'2.259"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").Value = "'0.01734"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").Value = "'0.4779"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'0.8966"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "'109.89"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").Value = "'5.929"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'7.681"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").Value = "'0.4158"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'9.033"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").Value = "'0.05811"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("D50").Value = "'34.76"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'0.8924"
$ws.Range("E51").Value = "  +0.67%  "
